$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure ambiguous numeric-looking price strings stay text (matches original inlineStr formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply cell value updates
$ws.Range("D2").Value = "27.683.46"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.844.16"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "313.09"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4283"
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("D8").Value = "0.3638"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.07317"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "0.8771"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "20.70"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").Value = "1.870.44"
$ws.Range("E12").Value = "  +2.56%  "
$ws.Range("D13").Value = "5.345"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "6.516"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").Value = "0.06937"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "0.000008991"
$ws.Range("E18").Value = "  +1.93%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "15.39"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "27.733.62"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "4.979"
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("D24").Value = "2.108.04"
$ws.Range("E24").Value = "  +1.40%  "
$ws.Range("D25").Value = "1.991"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "155.84"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("D27").Value = "18.56"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").Value = "119.67"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "1.883"
$ws.Range("E30").Value = "  +3.62%  "
$ws.Range("D31").Value = "0.08884"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "0.7530"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").Value = "4.527"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").Value = "2.956"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "1.117"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.05430"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "1.105"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "0.01934"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "2.829"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").Value = "0.1665"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "0.5067"
$ws.Range("E42").Value = "  +0.31%  "
$ws.Range("D43").Value = "6.596"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").Value = "0.06545"
$ws.Range("D46").Value = "10.37"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "105.86"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "0.4645"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("D50").Value = "1.634"
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "64.46"
$ws.Range("E51").Value = "  +0.07%  "
